$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.370.43"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "'2.646.69"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.85%  "
$ws.Range("D5").Value = "'598.24"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").Value = "'154.88"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.85%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D9").Value = "'2.645.56"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("E10").Value = "  +7.60%  "
$ws.Range("D12").Value = "'5.27"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("E13").Value = "  +2.33%  "
$ws.Range("D14").Value = "'28.06"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.89%  "
$ws.Range("E15").Value = "  +2.63%  "
$ws.Range("D16").Value = "'3.139.15"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.20%  "
$ws.Range("D17").Value = "'68.271.95"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").Value = "'2.645.01"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("D19").Value = "'11.37"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("D20").Value = "'364.16"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.79%  "
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("E22").Value = "  +3.25%  "
$ws.Range("E23").Value = "  +2.36%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").Value = "'74.85"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.92%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "'9.82"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.99%  "
$ws.Range("E28").Value = "  +2.49%  "
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("D31").Value = "'571.86"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.87%  "
$ws.Range("E32").Value = "  +2.64%  "
$ws.Range("E33").Value = "  +2.01%  "
$ws.Range("D34").Value = "'1.88"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.99%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("E37").Value = "  +5.47%  "
$ws.Range("D38").Value = "'161.15"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.19%  "
$ws.Range("E39").Value = "  +1.20%  "
$ws.Range("E40").Value = "  +2.08%  "
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("E42").Value = "  +1.00%  "
$ws.Range("D43").Value = "'0.0₆0338"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.24%  "
$ws.Range("E44").Value = "  +0.91%  "
$ws.Range("E45").Value = "  +2.02%  "
$ws.Range("D46").Value = "'40.63"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.05%  "
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("D48").Value = "'156.59"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.75%  "
$ws.Range("E49").Value = "  +2.07%  "
$ws.Range("E50").Value = "  +0.88%  "
$ws.Range("D51").Value = "'21.96"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.38%  "
